$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1743772241992882
$ws.Range("C2").Value = 0.5871886120996441
$ws.Range("J2").Value = 0.01067615658362989
$ws.Range("P2").Value = 0.1423487544483986
$ws.Range("S2").Value = 0.08540925266903915
$ws.Range("B3").Value = 0.02366863905325444
$ws.Range("C3").Value = 0.02366863905325444
$ws.Range("J3").Value = 0.05325443786982249
$ws.Range("P3").Value = 0.7218934911242604
$ws.Range("S3").Value = 0.1775147928994083
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.7096774193548387
$ws.Range("S4").Value = 0.2258064516129032
$ws.Range("B6").Value = 0.06481481481481481
$ws.Range("D6").Value = 0.02314814814814815
$ws.Range("F6").Value = 0.04166666666666666
$ws.Range("J6").Value = 0.3055555555555556
$ws.Range("O6").Value = 0.01388888888888889
$ws.Range("Q6").Value = 0.1296296296296296
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.337962962962963
$ws.Range("B7").Value = 0.095
$ws.Range("D7").Value = 0.02
$ws.Range("E7").Value = 0.005
$ws.Range("F7").Value = 0.045
$ws.Range("J7").Value = 0.115
$ws.Range("O7").Value = 0.025
$ws.Range("Q7").Value = 0.165
$ws.Range("R7").Value = 0.055
$ws.Range("S7").Value = 0.475
$ws.Range("B8").Value = 0.08426966292134831
$ws.Range("D8").Value = 0.01685393258426966
$ws.Range("F8").Value = 0.0599250936329588
$ws.Range("J8").Value = 0.1123595505617977
$ws.Range("O8").Value = 0.02434456928838951
$ws.Range("Q8").Value = 0.1629213483146068
$ws.Range("R8").Value = 0.08801498127340825
$ws.Range("S8").Value = 0.451310861423221
$ws.Range("B9").Value = 0.08666666666666667
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.1066666666666667
$ws.Range("O9").Value = 0.006666666666666667
$ws.Range("Q9").Value = 0.18
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.46
$ws.Range("B10").Value = 0.1084634346754314
$ws.Range("D10").Value = 0.01232539030402629
$ws.Range("F10").Value = 0.06655710764174198
$ws.Range("J10").Value = 0.1150369761709121
$ws.Range("O10").Value = 0.0180772391125719
$ws.Range("Q10").Value = 0.2013147082990961
$ws.Range("R10").Value = 0.07641741988496302
$ws.Range("S10").Value = 0.4018077239112572
$ws.Range("G11").Value = 0.1557632398753894
$ws.Range("J11").Value = 0.08411214953271028
$ws.Range("K11").Value = 0.1931464174454829
$ws.Range("L11").Value = 0.5607476635514018
$ws.Range("S11").Value = 0.006230529595015576
$ws.Range("G12").Value = 0.6989247311827957
$ws.Range("J12").Value = 0.2311827956989247
$ws.Range("K12").Value = 0.01075268817204301
$ws.Range("L12").Value = 0.02150537634408602
$ws.Range("S12").Value = 0.03763440860215054
$ws.Range("G13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.07843137254901961
$ws.Range("F15").Value = 0.0273972602739726
$ws.Range("H15").Value = 0.2146118721461187
$ws.Range("I15").Value = 0.0639269406392694
$ws.Range("J15").Value = 0.2694063926940639
$ws.Range("K15").Value = 0.0821917808219178
$ws.Range("M15").Value = 0.0273972602739726
$ws.Range("O15").Value = 0.1004566210045662
$ws.Range("S15").Value = 0.2146118721461187
$ws.Range("F16").Value = 0.02793296089385475
$ws.Range("H16").Value = 0.1955307262569832
$ws.Range("I16").Value = 0.08379888268156424
$ws.Range("J16").Value = 0.3798882681564246
$ws.Range("K16").Value = 0.1229050279329609
$ws.Range("M16").Value = 0.0111731843575419
$ws.Range("O16").Value = 0.0670391061452514
$ws.Range("S16").Value = 0.111731843575419
$ws.Range("F17").Value = 0.01900237529691211
$ws.Range("H17").Value = 0.2114014251781473
$ws.Range("I17").Value = 0.05463182897862233
$ws.Range("J17").Value = 0.3895486935866984
$ws.Range("K17").Value = 0.1045130641330166
$ws.Range("M17").Value = 0.02137767220902613
$ws.Range("O17").Value = 0.07363420427553444
$ws.Range("S17").Value = 0.1258907363420428
$ws.Range("F18").Value = 0.01657458563535912
$ws.Range("H18").Value = 0.2430939226519337
$ws.Range("I18").Value = 0.05524861878453038
$ws.Range("J18").Value = 0.3867403314917127
$ws.Range("K18").Value = 0.1160220994475138
$ws.Range("M18").Value = 0.01657458563535912
$ws.Range("O18").Value = 0.04972375690607735
$ws.Range("S18").Value = 0.1160220994475138
$ws.Range("F19").Value = 0.01656626506024097
$ws.Range("H19").Value = 0.2417168674698795
$ws.Range("I19").Value = 0.06852409638554217
$ws.Range("J19").Value = 0.3509036144578313
$ws.Range("K19").Value = 0.1129518072289157
$ws.Range("M19").Value = 0.02635542168674699
$ws.Range("N19").Value = 0.0007530120481927711
$ws.Range("O19").Value = 0.05798192771084337
$ws.Range("S19").Value = 0.1242469879518072
